# Sending testdata through excel
# - Adds a new "validLogin" worksheet after the existing "invalidLogin" sheet,
#   populates it with header + 2 data rows, and makes it the active sheet.
# - Clears the old "tabSelected"/activeCell selection on "invalidLogin" and
#   selects A1:C1 there instead.

$wb = $excel.ActiveWorkbook
$invalidLogin = $wb.Worksheets.Item(1)

# Update the selection on the existing sheet (no longer the active tab).
[void]$invalidLogin.Range("A1:C1").Select()

# Add the new sheet right after "invalidLogin" -- this also makes it active.
$validLogin = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $invalidLogin)
$validLogin.Name = "validLogin"

# Header row
$validLogin.Range("A1").Value = "username "
$validLogin.Range("B1").Value = "password "
$validLogin.Range("C1").Value = "expectedHeader"
$validLogin.Range("C1").WrapText = $true
$validLogin.Rows.Item(1).RowHeight = 30

# Data rows
$validLogin.Range("A2").Value = "Admin"
$validLogin.Range("B2").Value = "admin123"
$validLogin.Range("C2").Value = "Dashboard"

$validLogin.Range("A3").Value = "Admin"
$validLogin.Range("B3").Value = "admin123"
$validLogin.Range("C3").Value = "Dashboard"

# Match the saved selection/active-cell on the new sheet.
[void]$validLogin.Range("E4").Select()
